$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF), matching style of existing header cells
$ws.Range("I1").Value = "I0"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160
$ws.Range("I1").Borders.LineStyle = 1

$ws.Range("J1").Value = "IF"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1

# Data values for I (I0) and J (IF) columns, rows 2-70
$iVals = @{
    2 = 9
    3 = 9
    4 = 5
    5 = 8
    6 = 9
    7 = 6
    8 = 8
    9 = 1
    10 = 5
    11 = 7
    12 = 7
    13 = 10
    14 = 8
    15 = 5
    16 = 8
    17 = 7
    18 = 6
    19 = 9
    20 = 7
    21 = 7
    22 = 6
    23 = 7
    24 = 8
    25 = 6
    26 = 6
    27 = 6
    28 = 7
    29 = 6
    30 = 7
    31 = 7
    32 = 6
    33 = 6
    34 = 8
    35 = 6
    36 = 7
    37 = 10
    38 = 7
    39 = 7
    40 = 10
    41 = 8
    42 = 6
    43 = 6
    44 = 7
    45 = 7
    46 = 5
    47 = 1
    48 = 7
    49 = 1
    50 = 5
    51 = 6
    52 = 7
    53 = 1
    54 = 7
    55 = 1
    56 = 7
    57 = 8
    58 = 1
    59 = 1
    60 = 1
    61 = 1
    62 = 1
    63 = 1
    64 = 8
    65 = 1
    66 = 1
    67 = 1
    68 = 4
    69 = 4
    70 = 4
}
$jVals = @{
    2 = 9
    3 = 9
    4 = 5
    5 = 8
    6 = 9
    7 = 6
    8 = 8
    9 = 1
    10 = 5
    11 = 7
    12 = 7
    13 = 10
    14 = 8
    15 = 5
    16 = 8
    17 = 7
    18 = 7
    19 = 9
    20 = 7
    21 = 7
    22 = 7
    23 = 8
    24 = 8
    25 = 7
    26 = 6
    27 = 6
    28 = 7
    29 = 7
    30 = 8
    31 = 8
    32 = 7
    33 = 8
    34 = 8
    35 = 6
    36 = 7
    37 = 10
    38 = 8
    39 = 7
    40 = 11
    41 = 8
    42 = 7
    43 = 6
    44 = 7
    45 = 7
    46 = 6
    47 = 3
    48 = 7
    49 = 2
    50 = 6
    51 = 6
    52 = 8
    53 = 2
    54 = 7
    55 = 3
    56 = 7
    57 = 9
    58 = 4
    59 = 6
    60 = 6
    61 = 3
    62 = 5
    63 = 4
    64 = 8
    65 = 3
    66 = 3
    67 = 2
    68 = 4
    69 = 4
    70 = 4
}

foreach ($r in 2..70) {
    $ws.Cells.Item($r, 9).Value = $iVals[$r]
    $ws.Cells.Item($r, 10).Value = $jVals[$r]
}
